$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 60
$ws.Range("B60").Value = 6879734
$ws.Range("E60").Value = 'FC Petrzalka'
$ws.Range("F60").Value = 'OFK Malzenice'
$ws.Range("G60").Value = 2
$ws.Range("H60").Value = 2
$ws.Range("I60").Value = 'D'
$ws.Range("J60").Value = 1.25
$ws.Range("K60").Value = 5
$ws.Range("L60").Value = 8.5
$ws.Range("M60").Value = 1.363
$ws.Range("N60").Value = 5
$ws.Range("O60").Value = 7.5
$ws.Range("P60").Value = -1.25
$ws.Range("Q60").Value = 1.75
$ws.Range("R60").Value = 1.95
$ws.Range("S60").Value = 3
$ws.Range("T60").Value = 2
$ws.Range("U60").Value = 1.8
$ws.Range("V60").Value = -1
$ws.Range("W60").Value = 4
$ws.Range("X60").Value = -1
$ws.Range("Y60").Value = -1
$ws.Range("Z60").Value = 0.95
$ws.Range("AA60").Value = 1
$ws.Range("AB60").Value = -1

# Row 61
$ws.Range("B61").Value = 6883797
$ws.Range("E61").Value = 'STK Samorin'
$ws.Range("F61").Value = 'FK Tatran Liptovsky Mikulas'
$ws.Range("G61").Value = 3
$ws.Range("H61").Value = 1
$ws.Range("I61").Value = 'H'
$ws.Range("J61").Value = 2.625
$ws.Range("K61").Value = 3.5
$ws.Range("L61").Value = 2.25
$ws.Range("M61").Value = 4
$ws.Range("N61").Value = 3.8
$ws.Range("O61").Value = 1.8
$ws.Range("P61").Value = 0.5
$ws.Range("Q61").Value = 2
$ws.Range("R61").Value = 1.8
$ws.Range("S61").Value = 3
$ws.Range("T61").Value = 1.9
$ws.Range("U61").Value = 1.9
$ws.Range("V61").Value = 3
$ws.Range("W61").Value = -1
$ws.Range("X61").Value = -1
$ws.Range("Y61").Value = 1
$ws.Range("Z61").Value = -1
$ws.Range("AA61").Value = 0.8999999999999999
$ws.Range("AB61").Value = -1

# Row 139
$ws.Range("B139").Value = 6884078
$ws.Range("E139").Value = 'Spisska Nova Ves'
$ws.Range("F139").Value = 'FK Pohronie'
$ws.Range("G139").Value = 0
$ws.Range("H139").Value = 1
$ws.Range("I139").Value = 'A'
$ws.Range("J139").Value = 4.333
$ws.Range("K139").Value = 3.75
$ws.Range("L139").Value = 1.615
$ws.Range("M139").Value = 5
$ws.Range("N139").Value = 4.2
$ws.Range("O139").Value = 1.6
$ws.Range("P139").Value = 0.75
$ws.Range("Q139").Value = 1.95
$ws.Range("R139").Value = 1.75
$ws.Range("S139").Value = 2.75
$ws.Range("T139").Value = 1.875
$ws.Range("U139").Value = 1.925
$ws.Range("V139").Value = -1
$ws.Range("W139").Value = -1
$ws.Range("X139").Value = 0.6000000000000001
$ws.Range("Y139").Value = -0.5
$ws.Range("Z139").Value = 0.375
$ws.Range("AA139").Value = -1
$ws.Range("AB139").Value = 0.925

# Row 141
$ws.Range("B141").Value = 7923546
$ws.Range("E141").Value = 'OFK Malzenice'
$ws.Range("F141").Value = 'Puchov'
$ws.Range("G141").Value = 0
$ws.Range("H141").Value = 2
$ws.Range("I141").Value = 'A'
$ws.Range("J141").Value = 2.5
$ws.Range("K141").Value = 3.4
$ws.Range("L141").Value = 2.4
$ws.Range("M141").Value = 2.875
$ws.Range("N141").Value = 3.5
$ws.Range("O141").Value = 2.375
$ws.Range("P141").Value = 0.25
$ws.Range("Q141").Value = 1.775
$ws.Range("R141").Value = 2.025
$ws.Range("S141").Value = 2.5
$ws.Range("T141").Value = 1.8
$ws.Range("U141").Value = 2
$ws.Range("V141").Value = -1
$ws.Range("W141").Value = -1
$ws.Range("X141").Value = 1.375
$ws.Range("Y141").Value = -1
$ws.Range("Z141").Value = 1.025
$ws.Range("AA141").Value = -1
$ws.Range("AB141").Value = 1

# Row 142
$ws.Range("B142").Value = 6883465
$ws.Range("E142").Value = 'FK Humenne'
$ws.Range("F142").Value = 'Slovan Bratislava B'
$ws.Range("G142").Value = 2
$ws.Range("H142").Value = 0
$ws.Range("I142").Value = 'H'
$ws.Range("J142").Value = 1.909
$ws.Range("K142").Value = 3.3
$ws.Range("L142").Value = 3.5
$ws.Range("M142").Value = 1.533
$ws.Range("N142").Value = 4.2
$ws.Range("O142").Value = 6
$ws.Range("P142").Value = -1
$ws.Range("Q142").Value = 1.875
$ws.Range("R142").Value = 1.925
$ws.Range("S142").Value = 2.75
$ws.Range("T142").Value = 1.95
$ws.Range("U142").Value = 1.85
$ws.Range("V142").Value = 0.5329999999999999
$ws.Range("W142").Value = -1
$ws.Range("X142").Value = -1
$ws.Range("Y142").Value = 0.875
$ws.Range("Z142").Value = -1
$ws.Range("AA142").Value = -1
$ws.Range("AB142").Value = 0.8500000000000001

# Row 179
$ws.Range("B179").Value = 6883808
$ws.Range("E179").Value = 'FK Tatran Liptovsky Mikulas'
$ws.Range("F179").Value = 'MSK Povazska Bystrica'
$ws.Range("G179").Value = 1
$ws.Range("H179").Value = 2
$ws.Range("I179").Value = 'A'
$ws.Range("J179").Value = 2
$ws.Range("K179").Value = 3.5
$ws.Range("L179").Value = 3.05
$ws.Range("M179").Value = 1.65
$ws.Range("N179").Value = 4.2
$ws.Range("O179").Value = 5
$ws.Range("P179").Value = -0.75
$ws.Range("Q179").Value = 1.825
$ws.Range("R179").Value = 1.975
$ws.Range("S179").Value = 3
$ws.Range("T179").Value = 1.95
$ws.Range("U179").Value = 1.85
$ws.Range("V179").Value = -1
$ws.Range("W179").Value = -1
$ws.Range("X179").Value = 4
$ws.Range("Y179").Value = -1
$ws.Range("Z179").Value = 0.9750000000000001
$ws.Range("AA179").Value = 0
$ws.Range("AB179").Value = 0

# Row 181
$ws.Range("B181").Value = 6883475
$ws.Range("E181").Value = 'Slavoj Trebisov'
$ws.Range("F181").Value = 'FK Humenne'
$ws.Range("G181").Value = 1
$ws.Range("H181").Value = 1
$ws.Range("I181").Value = 'D'
$ws.Range("J181").Value = 3
$ws.Range("K181").Value = 3.4
$ws.Range("L181").Value = 2.05
$ws.Range("M181").Value = 5.75
$ws.Range("N181").Value = 4.2
$ws.Range("O181").Value = 1.6
$ws.Range("P181").Value = 1
$ws.Range("Q181").Value = 1.75
$ws.Range("R181").Value = 1.95
$ws.Range("S181").Value = 2.5
$ws.Range("T181").Value = 1.925
$ws.Range("U181").Value = 1.875
$ws.Range("V181").Value = -1
$ws.Range("W181").Value = 3.2
$ws.Range("X181").Value = -1
$ws.Range("Y181").Value = 0.75
$ws.Range("Z181").Value = -1
$ws.Range("AA181").Value = -1
$ws.Range("AB181").Value = 0.875
